$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column D header
$ws.Range("D1").Value = "store_code"

# Add store_code values for existing rows that have a built_template_code (rows 2 and 3)
$ws.Range("D2").Value = "shopee"
$ws.Range("D3").Value = "shopee"

# Remove the autofilter that covered A1:C8
$ws.AutoFilterMode = $false

# Update the active selection to D3
[void]$ws.Range("D3").Select()
